$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new CRUD "delete" test case rows to the Policies sheet
$ws.Range("A7").Value = "Verify that user is able to delete policy from Parent Company"
$ws.Range("B7").Value = "Verify that user is able to delete policy from Parent Company"
$ws.Range("C7").Value = "yes"

$ws.Range("A8").Value = "Verify that user is able to delete policy from Group Company"
$ws.Range("B8").Value = "Verify that user is able to delete policy from Group Company"
$ws.Range("C8").Value = "no"

# Match the author's final cursor/selection position
$ws.Range("A8").Select() | Out-Null
